$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2239382239382239
$ws.Range("C2").Value = 0.5057915057915058
$ws.Range("J2").Value = 0.01158301158301158
$ws.Range("P2").Value = 0.1544401544401544
$ws.Range("S2").Value = 0.1042471042471042
$ws.Range("C3").Value = 0.02238805970149254
$ws.Range("J3").Value = 0.02985074626865672
$ws.Range("P3").Value = 0.8208955223880597
$ws.Range("S3").Value = 0.1268656716417911
$ws.Range("J4").Value = 0.02631578947368421
$ws.Range("P4").Value = 0.7894736842105263
$ws.Range("S4").Value = 0.1842105263157895
$ws.Range("B6").Value = 0.07368421052631578
$ws.Range("D6").Value = 0.005263157894736842
$ws.Range("F6").Value = 0.03684210526315789
$ws.Range("J6").Value = 0.3368421052631579
$ws.Range("O6").Value = 0.01052631578947368
$ws.Range("Q6").Value = 0.1684210526315789
$ws.Range("R6").Value = 0.03157894736842105
$ws.Range("S6").Value = 0.3368421052631579
$ws.Range("B7").Value = 0.1282051282051282
$ws.Range("D7").Value = 0.01923076923076923
$ws.Range("F7").Value = 0.04487179487179487
$ws.Range("J7").Value = 0.1346153846153846
$ws.Range("O7").Value = 0.00641025641025641
$ws.Range("Q7").Value = 0.2371794871794872
$ws.Range("R7").Value = 0.03205128205128205
$ws.Range("S7").Value = 0.3974358974358974
$ws.Range("B8").Value = 0.07225433526011561
$ws.Range("D8").Value = 0.01734104046242774
$ws.Range("E8").Value = 0.002890173410404624
$ws.Range("F8").Value = 0.05202312138728324
$ws.Range("J8").Value = 0.1271676300578035
$ws.Range("O8").Value = 0.02023121387283237
$ws.Range("Q8").Value = 0.1878612716763006
$ws.Range("R8").Value = 0.07803468208092486
$ws.Range("S8").Value = 0.4421965317919075
$ws.Range("B9").Value = 0.1016949152542373
$ws.Range("D9").Value = 0.01129943502824859
$ws.Range("F9").Value = 0.06214689265536723
$ws.Range("J9").Value = 0.1638418079096045
$ws.Range("O9").Value = 0.02824858757062147
$ws.Range("Q9").Value = 0.1977401129943503
$ws.Range("R9").Value = 0.05084745762711865
$ws.Range("S9").Value = 0.384180790960452
$ws.Range("B10").Value = 0.09975669099756691
$ws.Range("D10").Value = 0.0218978102189781
$ws.Range("F10").Value = 0.0575831305758313
$ws.Range("J10").Value = 0.129764801297648
$ws.Range("O10").Value = 0.0145985401459854
$ws.Range("Q10").Value = 0.2327656123276561
$ws.Range("R10").Value = 0.06731549067315491
$ws.Range("S10").Value = 0.3763179237631792
$ws.Range("F11").Value = 0.007518796992481203
$ws.Range("G11").Value = 0.1691729323308271
$ws.Range("J11").Value = 0.07518796992481203
$ws.Range("K11").Value = 0.2330827067669173
$ws.Range("L11").Value = 0.4962406015037594
$ws.Range("S11").Value = 0.01879699248120301
$ws.Range("G12").Value = 0.6546762589928058
$ws.Range("J12").Value = 0.2302158273381295
$ws.Range("K12").Value = 0.02158273381294964
$ws.Range("L12").Value = 0.04316546762589928
$ws.Range("S12").Value = 0.05035971223021583
$ws.Range("G13").Value = 0.6153846153846154
$ws.Range("J13").Value = 0.3076923076923077
$ws.Range("S13").Value = 0.07692307692307693
$ws.Range("F15").Value = 0.02764976958525346
$ws.Range("H15").Value = 0.1105990783410138
$ws.Range("I15").Value = 0.07834101382488479
$ws.Range("J15").Value = 0.3963133640552995
$ws.Range("K15").Value = 0.03225806451612903
$ws.Range("M15").Value = 0.009216589861751152
$ws.Range("O15").Value = 0.06451612903225806
$ws.Range("S15").Value = 0.2811059907834101
$ws.Range("F16").Value = 0.02339181286549707
$ws.Range("H16").Value = 0.1637426900584795
$ws.Range("I16").Value = 0.1169590643274854
$ws.Range("J16").Value = 0.3801169590643275
$ws.Range("K16").Value = 0.09941520467836257
$ws.Range("M16").Value = 0.01169590643274854
$ws.Range("O16").Value = 0.08187134502923976
$ws.Range("S16").Value = 0.1228070175438596
$ws.Range("F17").Value = 0.02838427947598253
$ws.Range("H17").Value = 0.1397379912663755
$ws.Range("I17").Value = 0.08078602620087336
$ws.Range("J17").Value = 0.4432314410480349
$ws.Range("K17").Value = 0.08733624454148471
$ws.Range("M17").Value = 0.01091703056768559
$ws.Range("O17").Value = 0.0851528384279476
$ws.Range("S17").Value = 0.1244541484716157
$ws.Range("F18").Value = 0.0310077519379845
$ws.Range("H18").Value = 0.1395348837209302
$ws.Range("I18").Value = 0.1007751937984496
$ws.Range("J18").Value = 0.4108527131782946
$ws.Range("K18").Value = 0.06201550387596899
$ws.Range("M18").Value = 0.0310077519379845
$ws.Range("O18").Value = 0.07751937984496124
$ws.Range("S18").Value = 0.1472868217054264
$ws.Range("F19").Value = 0.02200704225352113
$ws.Range("H19").Value = 0.1936619718309859
$ws.Range("I19").Value = 0.07746478873239436
$ws.Range("J19").Value = 0.3961267605633803
$ws.Range("K19").Value = 0.1073943661971831
$ws.Range("M19").Value = 0.02288732394366197
$ws.Range("O19").Value = 0.07130281690140845
$ws.Range("S19").Value = 0.1091549295774648
